$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new weekly rows before the old row 250 (rows 250-259 shift down to 253-262)
$ws.Range("A250:A252").EntireRow.Insert()

# Row 250
$ws.Cells.Item(250, 1).Value = 6
$ws.Cells.Item(250, 2).Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Cells.Item(250, 3).Value = 'Metropolitana'
$ws.Cells.Item(250, 4).Value = 44516
$ws.Cells.Item(250, 5).Value = 13
$ws.Cells.Item(250, 6).Value = 100112032
$ws.Cells.Item(250, 7).Value = 'Zapallo italiano'
$ws.Cells.Item(250, 8).Value = 'Sin especificar'
$ws.Cells.Item(250, 9).Value = 'Primera'
$ws.Cells.Item(250, 10).Value = 400
$ws.Cells.Item(250, 11).Value = 5000
$ws.Cells.Item(250, 12).Value = 6000
$ws.Cells.Item(250, 13).Value = 5425
$ws.Cells.Item(250, 14).Value = '$/caja 50 unidades'
$ws.Cells.Item(250, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(250, 16).Value = 108
$ws.Cells.Item(250, 17).Value = 50
$ws.Cells.Item(250, 18).Value = 'Hortaliza'

# Row 251
$ws.Cells.Item(251, 1).Value = 6
$ws.Cells.Item(251, 2).Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Cells.Item(251, 3).Value = 'Metropolitana'
$ws.Cells.Item(251, 4).Value = 44516
$ws.Cells.Item(251, 5).Value = 13
$ws.Cells.Item(251, 6).Value = 100112032
$ws.Cells.Item(251, 7).Value = 'Zapallo italiano'
$ws.Cells.Item(251, 8).Value = 'Sin especificar'
$ws.Cells.Item(251, 9).Value = 'Primera'
$ws.Cells.Item(251, 10).Value = 1030
$ws.Cells.Item(251, 11).Value = 5000
$ws.Cells.Item(251, 12).Value = 6000
$ws.Cells.Item(251, 13).Value = 5544
$ws.Cells.Item(251, 14).Value = '$/caja 50 unidades'
$ws.Cells.Item(251, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(251, 16).Value = 111
$ws.Cells.Item(251, 17).Value = 50
$ws.Cells.Item(251, 18).Value = 'Hortaliza'

# Row 252
$ws.Cells.Item(252, 1).Value = 6
$ws.Cells.Item(252, 2).Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Cells.Item(252, 3).Value = 'Metropolitana'
$ws.Cells.Item(252, 4).Value = 44516
$ws.Cells.Item(252, 5).Value = 13
$ws.Cells.Item(252, 6).Value = 100112032
$ws.Cells.Item(252, 7).Value = 'Zapallo italiano'
$ws.Cells.Item(252, 8).Value = 'Sin especificar'
$ws.Cells.Item(252, 9).Value = 'Segunda'
$ws.Cells.Item(252, 10).Value = 400
$ws.Cells.Item(252, 11).Value = 3000
$ws.Cells.Item(252, 12).Value = 4000
$ws.Cells.Item(252, 13).Value = 3575
$ws.Cells.Item(252, 14).Value = '$/caja 80 unidades'
$ws.Cells.Item(252, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(252, 16).Value = 45
$ws.Cells.Item(252, 17).Value = 80
$ws.Cells.Item(252, 18).Value = 'Hortaliza'
